$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # planificación
$ws2 = $wb.Worksheets.Item(2)   # Restricciones de Usuario
$ws3 = $wb.Worksheets.Item(3)   # Calendario fechas

# --- New rows of restriction data in "Restricciones de Usuario" ---
$ws2.Range("A6").Value = "TD"
$ws2.Range("B6").Value = 13
$ws2.Range("C6").Value = 38
$ws2.Range("D6").Value = 1

$ws2.Range("D5").Value = 3

$ws2.Range("A7").Value = "SD"
$ws2.Range("B7").Value = 46
$ws2.Range("C7").Value = 47

$ws2.Range("A8").Value = "SD"
$ws2.Range("B8").Value = 34
$ws2.Range("C8").Value = 46

$ws2.Range("A9").Value = "SD"
$ws2.Range("B9").Value = 47
$ws2.Range("C9").Value = 3

$ws2.Range("A10").Value = "DB"
$ws2.Range("B10").Value = 3
$ws2.Range("C10").Value = 44361

$ws2.Range("A11").Value = "DB"
$ws2.Range("B11").Value = 3
$ws2.Range("C11").Value = 44362

$ws2.Range("A12").Value = "DB"
$ws2.Range("B12").Value = 3
$ws2.Range("C12").Value = 44363

$ws2.Range("A13").Value = "DB"
$ws2.Range("B13").Value = 3
$ws2.Range("C13").Value = 44365

$ws2.Range("A14").Value = "DB"
$ws2.Range("B14").Value = 3
$ws2.Range("C14").Value = 44368

$ws2.Range("A15").Value = "DB"
$ws2.Range("B15").Value = 3
$ws2.Range("C15").Value = 44369

$ws2.Range("A16").Value = "DB"
$ws2.Range("B16").Value = 3
$ws2.Range("C16").Value = 44370

$ws2.Range("A17").Value = "DB"
$ws2.Range("B17").Value = 3
$ws2.Range("C17").Value = 44371

$ws2.Range("A18").Value = "DB"
$ws2.Range("B18").Value = 3
$ws2.Range("C18").Value = 44372

$ws2.Range("A19").Value = "DB"
$ws2.Range("B19").Value = 3
$ws2.Range("C19").Value = 44375

$ws2.Range("A20").Value = "DB"
$ws2.Range("B20").Value = 3
$ws2.Range("C20").Value = 44376

$ws2.Range("A21").Value = "DB"
$ws2.Range("B21").Value = 3
$ws2.Range("C21").Value = 44377

$ws2.Range("A22").Value = "DB"
$ws2.Range("B22").Value = 3
$ws2.Range("C22").Value = 44378

$ws2.Range("A23").Value = "DB"
$ws2.Range("B23").Value = 3
$ws2.Range("C23").Value = 44379

$ws2.Range("A24").Value = "DB"
$ws2.Range("B24").Value = 3
$ws2.Range("C24").Value = 44382

$ws2.Range("A25").Value = "DB"
$ws2.Range("B25").Value = 3
$ws2.Range("C25").Value = 44383

# --- Copy number formats for the styled cells ---
$ws2.Range("C4").Copy()
$ws2.Range("C6:C9").PasteSpecial(-4122)

$ws3.Range("A1").Copy()
$ws2.Range("C10").PasteSpecial(-4122)

$ws3.Range("A2").Copy()
$ws2.Range("C11:C25").PasteSpecial(-4122)

$ws2.Range("E16").Font.Bold = $true

# --- Conditional formatting range updates ---
$ws2.Range("C10:C25").FormatConditions.Item(1).Delete()

# --- Activate the "Restricciones de Usuario" sheet and select E16 ---
$ws2.Activate()
$ws2.Range("E16").Select()
